# Milestone 1 folder update for submission:
#  - assign team members to backlog tasks in column F (rows 5-39)
#  - row 32 grows to a 2-line height once its F cell gets text
#  - leave the final selection on F5:F39 (matches the on-screen state
#    after filling the column) and activate the Product Backlog sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")
$ws.Activate()

# "Assigned to" values (column F) for rows 5-39, grouped by contiguous runs
$ws.Range("F5:F16").Value = "Luke, Ash"
$ws.Range("F17:F22").Value = "Joseph, Abdul"
$ws.Range("F23:F26").Value = "Amy"
$ws.Range("F27:F32").Value = "Joseph, Abdul"
$ws.Range("F33:F38").Value = "Amy"
$ws.Range("F39").Value = "Amy "

# Row 32 had no custom height before; once F32 carries text under the
# wrap-text style it needs two lines to display, so bump its height.
$ws.Rows.Item(32).RowHeight = 25.5

# Leave the sheet scrolled/selected the way it was when the edits were made
$ws.Range("F5:F39").Select() | Out-Null
